$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at X (column 24), shifting the existing
# "profit/expected payoff/drawdown url" columns (previously X:AF) one
# position to the right (now Y:AG).
$ws.Columns("X:X").Insert()

# The freshly inserted column inherits the formatting of its left
# neighbour (standard Excel insert behaviour) - clear it so the column
# stays blank/unformatted, same as the author's edit.
$ws.Range("X1").Clear()

# Update the active selection to match the edited sheet.
$ws.Range("AB3").Select()
